$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3802213333333334
$ws.Range("H2").Value = 1.140664
$ws.Range("I2").Value = 0.05154022338265814
$ws.Range("J2").Value = 0.05154022338265814
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05256533333333333
$ws.Range("N2").Value = 0.157696
$ws.Range("O2").Value = 0.7714955259952154
$ws.Range("P2").Value = 0.7714955259952153
$ws.Range("Q2").Value = 0.01998646112711111
$ws.Range("R2").Value = 0.179878150144
$ws.Range("S2").Value = 0.03976305174851474
$ws.Range("T2").Value = 0.03976305174851474
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3802213333333334
$ws.Range("H3").Value = 1.140664
$ws.Range("I3").Value = 0.05154022338265814
$ws.Range("J3").Value = 0.05154022338265814
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.015569
$ws.Range("N3").Value = 0.046707
$ws.Range("O3").Value = 0.2285044740047847
$ws.Range("P3").Value = 0.2285044740047847
$ws.Range("Q3").Value = 0.005919665938666667
$ws.Range("R3").Value = 0.05327699344800001
$ws.Range("S3").Value = 0.0117771716341434
$ws.Range("T3").Value = 0.0117771716341434
$ws.Range("G4").Value = 6.845175999999999
$ws.Range("H4").Value = 20.535528
$ws.Range("I4").Value = 0.9278856003177367
$ws.Range("J4").Value = 0.9278856003177366
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05256533333333333
$ws.Range("N4").Value = 0.157696
$ws.Range("O4").Value = 0.7714955259952154
$ws.Range("P4").Value = 0.7714955259952153
$ws.Range("Q4").Value = 0.3598189581653333
$ws.Range("R4").Value = 3.238370623488
$ws.Range("S4").Value = 0.7158595892805185
$ws.Range("T4").Value = 0.7158595892805183
$ws.Range("G5").Value = 6.845175999999999
$ws.Range("H5").Value = 20.535528
$ws.Range("I5").Value = 0.9278856003177367
$ws.Range("J5").Value = 0.9278856003177366
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.015569
$ws.Range("N5").Value = 0.046707
$ws.Range("O5").Value = 0.2285044740047847
$ws.Range("P5").Value = 0.2285044740047847
$ws.Range("Q5").Value = 0.106572545144
$ws.Range("R5").Value = 0.9591529062959999
$ws.Range("S5").Value = 0.2120260110372183
$ws.Range("T5").Value = 0.2120260110372182
$ws.Range("G6").Value = 0.1517793333333333
$ws.Range("H6").Value = 0.455338
$ws.Range("I6").Value = 0.02057417629960514
$ws.Range("J6").Value = 0.02057417629960513
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05256533333333333
$ws.Range("N6").Value = 0.157696
$ws.Range("O6").Value = 0.7714955259952154
$ws.Range("P6").Value = 0.7714955259952153
$ws.Range("Q6").Value = 0.007978331249777778
$ws.Range("R6").Value = 0.071804981248
$ws.Range("S6").Value = 0.01587288496618216
$ws.Range("T6").Value = 0.01587288496618216
$ws.Range("G7").Value = 0.1517793333333333
$ws.Range("H7").Value = 0.455338
$ws.Range("I7").Value = 0.02057417629960514
$ws.Range("J7").Value = 0.02057417629960513
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.015569
$ws.Range("N7").Value = 0.046707
$ws.Range("O7").Value = 0.2285044740047847
$ws.Range("P7").Value = 0.2285044740047847
$ws.Range("Q7").Value = 0.002363052440666667
$ws.Range("R7").Value = 0.021267471966
$ws.Range("S7").Value = 0.004701291333422979
$ws.Range("T7").Value = 0.004701291333422977
